$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 (pushes existing rows 9-51 down to 10-52).
$ws.Rows("9:9").Insert()

# Populate the new row 9 with a fresh weekly price observation. It mirrors
# the (now shifted-down) row 10 data except for the date and price columns,
# which carry this week's new figures.
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44561
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112021
$ws.Cells.Item(9, 7).Value = "Ají"
$ws.Cells.Item(9, 8).Value = "Americana (o)"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 60
$ws.Cells.Item(9, 11).Value = 20000
$ws.Cells.Item(9, 12).Value = 21000
$ws.Cells.Item(9, 13).Value = 20500
$ws.Cells.Item(9, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 1367
$ws.Cells.Item(9, 17).Value = 15
$ws.Cells.Item(9, 18).Value = "Hortaliza"
